$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3
$ws.Range("X2").Value = "free"
$ws.Range("AY2").Value = "free"
$ws.Range("BI2").Value = "free"
$ws.Range("BO2").Value = 0.33

$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$win.ScrollRow = 1
$ws.Range("Y7").Select()
